$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the "R10" rule row (cell E8)
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection on the sheet
$ws.Range("E8").Select()
